# Updated IPS AIP hipo turnover
# Target worksheet: "Rosemont Illinois" (corresponds to SOR Testing_Corp RBS.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rosemont Illinois")

# Professional Voluntary Turnover (cvd) - ytd value updated on rows 2-4
# (PY Actual, AOP, Commit/Forecast)
$ws.Range("E2").Value = 1.5385
$ws.Range("E3").Value = 1.5385
$ws.Range("E4").Value = 1.5385

# Row 4 (Commit/Forecast) monthly/quarterly turnover figures recalculated
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.256416666666667
$ws.Range("P4").Value = 0.256416666666667
$ws.Range("Q4").Value = 0.256416666666667
$ws.Range("R4").Value = 0.76925
$ws.Range("S4").Value = 0.256416666666667
$ws.Range("T4").Value = 0.256416666666667
$ws.Range("U4").Value = 0.256416666666667
$ws.Range("V4").Value = 0.76925
$ws.Range("W4").Value = 3.077

# Row 7 (Internal Fill Rate, Commit/Forecast): clear Apr/May values (now blank)
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
